$wb = $excel.ActiveWorkbook

# --- sheet ALC (hunk starting near old line 1483) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2619.75
$ws.Range("J17").Value = 2619.75
$ws.Range("L17").Value = 7859.25
$ws.Range("N17").Value = -8195.25

# --- sheet ALC (hunk starting near old line 3792) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 10000
$ws.Range("J64").Value = 10000
$ws.Range("L64").Value = 10000
$ws.Range("N64").Value = -10496

# --- sheet ALC (hunk starting near old line 3939) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 10000
$ws.Range("J67").Value = 10000
$ws.Range("L67").Value = 10000
$ws.Range("N67").Value = -11716

# --- sheet ALC (hunk starting near old line 4037) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 7593.7856
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 7593.7856
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 22781.3568
$ws.Range("M69").Value = ""
$ws.Range("N69").Value = -24529.3568

# --- sheet ALC (hunk starting near old line 4089) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 4884.1665
$ws.Range("I70").Value = 1818.25
$ws.Range("J70").Value = 6417.125
$ws.Range("K70").Value = 5454.75
$ws.Range("L70").Value = 19251.375
$ws.Range("M70").Value = -5184.75
$ws.Range("N70").Value = -19791.375

# --- sheet ALC (hunk starting near old line 4190) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H72").Value = 7593.7856
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 7593.7856
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 68344.0704
$ws.Range("M72").Value = ""
$ws.Range("N72").Value = -77080.0704

# --- sheet ALC (hunk starting near old line 4242) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 4884.1665
$ws.Range("I73").Value = 1818.25
$ws.Range("J73").Value = 6417.125
$ws.Range("K73").Value = 5454.75
$ws.Range("L73").Value = 19251.375
$ws.Range("M73").Value = -4518.75
$ws.Range("N73").Value = -21123.375

# --- sheet ALC (hunk starting near old line 6983) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1277
$ws.Range("J129").Value = 1599
$ws.Range("L129").Value = 4797
$ws.Range("N129").Value = -14797

# --- sheet ARM (hunk starting near old line 14070) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H133").Value = 89999
$ws.Range("J133").Value = 89999
$ws.Range("L133").Value = 89999
$ws.Range("N133").Value = -95059

# --- sheet BSM (hunk starting near old line 16842) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H49").Value = 63065
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 63065
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 63065
$ws.Range("M49").Value = ""
$ws.Range("N49").Value = -63543

# --- sheet BSM (hunk starting near old line 18732) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H88").Value = 28000
$ws.Range("J88").Value = 28000
$ws.Range("L88").Value = 28000
$ws.Range("N88").Value = -28812

# --- sheet BSM (hunk starting near old line 18879) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H91").Value = 28000
$ws.Range("J91").Value = 28000
$ws.Range("L91").Value = 28000
$ws.Range("N91").Value = -30808

# --- sheet CRP (hunk starting near old line 24198) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3140.4
$ws.Range("I58").Value = 1676.75
$ws.Range("K58").Value = 1676.75
$ws.Range("M58").Value = -1473.75

# --- sheet CRP (hunk starting near old line 25576) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 3908
$ws.Range("J86").Value = 3908
$ws.Range("L86").Value = 3908
$ws.Range("N86").Value = -6154

# --- sheet CRP (hunk starting near old line 25720) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 3908
$ws.Range("J89").Value = 3908
$ws.Range("L89").Value = 19540
$ws.Range("N89").Value = -30772

# --- sheet CRP (hunk starting near old line 27990) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 3140.4
$ws.Range("I136").Value = 1676.75
$ws.Range("K136").Value = 5030.25
$ws.Range("M136").Value = -2480.25

# --- sheet CUL (hunk starting near old line 29988) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1521.4286
$ws.Range("J34").Value = 999.2727
$ws.Range("L34").Value = 2997.8181
$ws.Range("N34").Value = -3165.8181

# --- sheet CUL (hunk starting near old line 30239) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 7036.6665
$ws.Range("J39").Value = 7036.6665
$ws.Range("L39").Value = 21109.9995
$ws.Range("N39").Value = -21697.9995

# --- sheet CUL (hunk starting near old line 31026) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 7777.6665
$ws.Range("J55").Value = 7777.6665
$ws.Range("L55").Value = 23332.9995
$ws.Range("N55").Value = -23686.9995

# --- sheet CUL (hunk starting near old line 34064) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 122.5
$ws.Range("I117").Value = 122.5
$ws.Range("J117").Value = 0
$ws.Range("K117").Value = 367.5
$ws.Range("L117").Value = 0
$ws.Range("M117").Value = 3074.5
$ws.Range("N117").Value = ""

# --- sheet CUL (hunk starting near old line 34453) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H125").Value = 7820
$ws.Range("J125").Value = 7980
$ws.Range("L125").Value = 23940
$ws.Range("N125").Value = -33780

# --- sheet GSM (hunk starting near old line 39145) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3841.5715
$ws.Range("I80").Value = 3819
$ws.Range("K80").Value = 3819
$ws.Range("M80").Value = -2821

# --- sheet GSM (hunk starting near old line 39243) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H82").Value = 50000
$ws.Range("J82").Value = 50000
$ws.Range("L82").Value = 50000
$ws.Range("N82").Value = -50766

# --- sheet GSM (hunk starting near old line 39289) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 3841.5715
$ws.Range("I83").Value = 3819
$ws.Range("K83").Value = 19095
$ws.Range("M83").Value = -14103

# --- sheet GSM (hunk starting near old line 39387) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H85").Value = 50000
$ws.Range("J85").Value = 50000
$ws.Range("L85").Value = 50000
$ws.Range("N85").Value = -52652

# --- sheet GSM (hunk starting near old line 41173) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4645.4
$ws.Range("I122").Value = 3075.6667
$ws.Range("K122").Value = 9227.000100000001
$ws.Range("M122").Value = -6777.000100000001

# --- sheet GSM (hunk starting near old line 41372) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 5270.6665
$ws.Range("I126").Value = 5270.6665
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 15811.9995
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -13341.9995
$ws.Range("N126").Value = ""

# --- sheet LTW (hunk starting near old line 44370) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 7357.143
$ws.Range("J46").Value = 8875
$ws.Range("L46").Value = 8875
$ws.Range("N46").Value = -9251

# --- sheet LTW (hunk starting near old line 45424) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 6100.1113
$ws.Range("I68").Value = 2450.5
$ws.Range("K68").Value = 2450.5
$ws.Range("M68").Value = -1701.5

# --- sheet LTW (hunk starting near old line 45574) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 6100.1113
$ws.Range("I71").Value = 2450.5
$ws.Range("K71").Value = 12252.5
$ws.Range("M71").Value = -8508.5

# --- sheet LTW (hunk starting near old line 46986) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 6495.5
$ws.Range("J100").Value = 9187.125
$ws.Range("L100").Value = 9187.125
$ws.Range("N100").Value = -10269.125

# --- sheet LTW (hunk starting near old line 48539) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 1342
$ws.Range("I132").Value = 1342
$ws.Range("K132").Value = 4026
$ws.Range("M132").Value = -1496

# --- sheet LTW (hunk starting near old line 48729) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1397.5
$ws.Range("I136").Value = 1397.5
$ws.Range("K136").Value = 4192.5
$ws.Range("M136").Value = -1642.5

# --- sheet WVR (hunk starting near old line 52910) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 5000
$ws.Range("I81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("M81").Value = ""

# --- sheet WVR (hunk starting near old line 53060) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 5000
$ws.Range("I84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("M84").Value = ""

# --- sheet WVR (hunk starting near old line 54178) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 859.3333
$ws.Range("I107").Value = 866.75
$ws.Range("J107").Value = 800
$ws.Range("K107").Value = 2600.25
$ws.Range("L107").Value = 2400
$ws.Range("M107").Value = -680.25
$ws.Range("N107").Value = -6240
